# Update NATMI LR-pairs data (Sertad1 -> Ar) with the refreshed TPM-based
# run: existing rows 2-16 get new edge-weight/specificity figures, the
# sending/target cluster pairing on rows 4-13 is re-ordered, and rows 16-21
# grow the table from 15 to 20 sender/target combinations (now including
# "Inflammatory-Mac" as a target cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per new record: Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then the 16 numeric NATMI statistics columns (E:T).
$data = @(
    @('ECs', 'Sertad1', 'Ar', 'ECs', 3, 1, 10.76654833333333, 32.299645, 0.1686652623327698, 0.175051555824113, 3, 1, 0.9764603333333333, 2.929381, 0.09919147117119548, 0.1117716869728403, 10.51310737441611, 94.61796636974499, 0.01673015550626306, 0.01956580770168144),
    @('ECs', 'Sertad1', 'Ar', 'FAPs', 3, 1, 10.76654833333333, 32.299645, 0.1686652623327698, 0.175051555824113, 3, 1, 5.511558666666666, 16.534676, 0.5598789770873293, 0.6308870816289637, 59.34046277666888, 534.0641649900199, 0.09443213454503721, 0.1104377651884843),
    @('ECs', 'Sertad1', 'Ar', 'Inflammatory-Mac', 3, 1, 10.76654833333333, 32.299645, 0.1686652623327698, 0.175051555824113, 1, 0.3333333333333333, 0.032202, 0.096606, 0.00327116590978248, 0.003686040017224872, 0.34670438943, 3.12033950487, 0.0005517320563074757, 0.000645247039845154),
    @('ECs', 'Sertad1', 'Ar', 'MuSCs', 3, 1, 10.76654833333333, 32.299645, 0.1686652623327698, 0.175051555824113, 2, 1, 3.3239755, 6.647951, 0.3376583858316928, 0.2536551913809712, 35.78774287956583, 214.726457277395, 0.05695124022516208, 0.04440273589410215),
    @('FAPs', 'Sertad1', 'Ar', 'ECs', 3, 1, 10.98286, 32.94858, 0.1720539309082887, 0.1785685319821705, 3, 1, 0.9764603333333333, 2.929381, 0.09919147117119548, 0.1117716869728403, 10.72432713655333, 96.51894422897999, 0.01706628252758038, 0.01995890605991078),
    @('FAPs', 'Sertad1', 'Ar', 'FAPs', 3, 1, 10.98286, 32.94858, 0.1720539309082887, 0.1785685319821705, 3, 1, 5.511558666666666, 16.534676, 0.5598789770873293, 0.6308870816289637, 60.53267721778666, 544.7940949600799, 0.09632937884078671, 0.1126565800129998),
    @('FAPs', 'Sertad1', 'Ar', 'Inflammatory-Mac', 3, 1, 10.98286, 32.94858, 0.1720539309082887, 0.1785685319821705, 1, 0.3333333333333333, 0.032202, 0.096606, 0.00327116590978248, 0.003686040017224872, 0.35367005772, 3.18303051948, 0.0005628169534312642, 0.0006582107547033798),
    @('FAPs', 'Sertad1', 'Ar', 'MuSCs', 3, 1, 10.98286, 32.94858, 0.1720539309082887, 0.1785685319821705, 2, 1, 3.3239755, 6.647951, 0.3376583858316928, 0.2536551913809712, 36.50675755993, 219.04054535958, 0.05809545258649037, 0.04529483515455653),
    @('Inflammatory-Mac', 'Sertad1', 'Ar', 'ECs', 3, 1, 23.02053966666667, 69.06161900000001, 0.3606323253943132, 0.3742872051281717, 3, 1, 0.9764603333333333, 2.929381, 0.09919147117119548, 0.1117716869728403, 22.47864383642656, 202.307794527839, 0.03577165090775121, 0.04183471232952528),
    @('Inflammatory-Mac', 'Sertad1', 'Ar', 'FAPs', 3, 1, 23.02053966666667, 69.06161900000001, 0.3606323253943132, 0.3742872051281717, 3, 1, 5.511558666666666, 16.534676, 0.5598789770873293, 0.6308870816289637, 126.8790549111604, 1141.911494200444, 0.201910457446393, 0.2361329625343735),
    @('Inflammatory-Mac', 'Sertad1', 'Ar', 'Inflammatory-Mac', 3, 1, 23.02053966666667, 69.06161900000001, 0.3606323253943132, 0.3742872051281717, 1, 0.3333333333333333, 0.032202, 0.096606, 0.00327116590978248, 0.003686040017224872, 0.7413074183460001, 6.671766765114, 0.00117968816879546, 0.001379637616037695),
    @('Inflammatory-Mac', 'Sertad1', 'Ar', 'MuSCs', 3, 1, 23.02053966666667, 69.06161900000001, 0.3606323253943132, 0.3742872051281717, 2, 1, 3.3239755, 6.647951, 0.3376583858316928, 0.2536551913809712, 76.51970984877818, 459.1182590926691, 0.1217705288713736, 0.09493989264823521),
    @('MuSCs', 'Sertad1', 'Ar', 'ECs', 2, 1, 6.986425000000001, 13.97285, 0.1094470733712295, 0.07572743080603386, 3, 1, 0.9764603333333333, 2.929381, 0.09919147117119548, 0.1117716869728403, 6.821966884308333, 40.93180130585, 0.01085621622307402, 0.008464182691309441),
    @('MuSCs', 'Sertad1', 'Ar', 'FAPs', 2, 1, 6.986425000000001, 13.97285, 0.1094470733712295, 0.07572743080603386, 3, 1, 5.511558666666666, 16.534676, 0.5598789770873293, 0.6308870816289637, 38.50609125776666, 231.0365475466, 0.06127711548428582, 0.04777545782047798),
    @('MuSCs', 'Sertad1', 'Ar', 'Inflammatory-Mac', 2, 1, 6.986425000000001, 13.97285, 0.1094470733712295, 0.07572743080603386, 1, 0.3333333333333333, 0.032202, 0.096606, 0.00327116590978248, 0.003686040017224872, 0.22497685785, 1.3498611471, 0.0003580195353374277, 0.0002791343403526683),
    @('MuSCs', 'Sertad1', 'Ar', 'MuSCs', 2, 1, 6.986425000000001, 13.97285, 0.1094470733712295, 0.07572743080603386, 2, 1, 3.3239755, 6.647951, 0.3376583858316928, 0.2536551913809712, 23.2227055325875, 92.89082213035, 0.03695572212853219, 0.01920865595389377),
    @('Resolving-Mac', 'Sertad1', 'Ar', 'ECs', 3, 1, 12.077449, 36.232347, 0.1892014079933989, 0.196365276259511, 3, 1, 0.9764603333333333, 2.929381, 0.09919147117119548, 0.1117716869728403, 11.79314987635633, 106.138348887207, 0.01876716600652682, 0.02194807819041337),
    @('Resolving-Mac', 'Sertad1', 'Ar', 'FAPs', 3, 1, 12.077449, 36.232347, 0.1892014079933989, 0.196365276259511, 3, 1, 5.511558666666666, 16.534676, 0.5598789770873293, 0.6308870816289637, 66.56556870717465, 599.0901183645718, 0.1059298907708266, 0.1238843160726281),
    @('Resolving-Mac', 'Sertad1', 'Ar', 'Inflammatory-Mac', 3, 1, 12.077449, 36.232347, 0.1892014079933989, 0.196365276259511, 1, 0.3333333333333333, 0.032202, 0.096606, 0.00327116590978248, 0.003686040017224872, 0.388918012698, 3.500262114282, 0.0006189091959108528, 0.0007238102662859746),
    @('Resolving-Mac', 'Sertad1', 'Ar', 'MuSCs', 3, 1, 12.077449, 36.232347, 0.1892014079933989, 0.196365276259511, 2, 1, 3.3239755, 6.647951, 0.3376583858316928, 0.2536551913809712, 40.1451445784995, 240.870867470997, 0.0638854420201346, 0.04980907173018354)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
